$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; existing rows 13-68 shift down to 14-69.
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with this week's data.
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(13, 3).Value = "Bíobío"
$ws.Cells.Item(13, 4).Value = 44670
$ws.Cells.Item(13, 5).Value = 8
$ws.Cells.Item(13, 6).Value = 100112012
$ws.Cells.Item(13, 7).Value = "Espinaca"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 170
$ws.Cells.Item(13, 11).Value = 5000
$ws.Cells.Item(13, 12).Value = 6000
$ws.Cells.Item(13, 13).Value = 5471
$ws.Cells.Item(13, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(13, 15).Value = "Región Metropolitana"
$ws.Cells.Item(13, 16).Value = 547
$ws.Cells.Item(13, 17).Value = 10
$ws.Cells.Item(13, 18).Value = "Hortaliza"
